$wb = $excel.ActiveWorkbook

# --- sheet "links": add new FATN / 葛章峰 mapping row ---
$links = $wb.Worksheets.Item("links")
$links.Range("A12").Value = "FATN"
$links.Range("B12").Value = "葛章峰"

# --- sheet "template": extend the dropdown list validations ---
$tmpl = $wb.Worksheets.Item("template")

$tmpl.Range("F3:F1048576").Validation.Delete()
$tmpl.Range("G3:G1048576").Validation.Delete()

$tmpl.Range("F3:F1048576").Validation.Add(3, 1, 1, '"DT部,VT部,SWT部,NPI部"')
$tmpl.Range("G3:G1048576").Validation.Add(3, 1, 1, '"SYD,HWD,MED,CSV,HWV,SSD,SCD,SWV,PSD,CUD,FWD,FATN"')
